# Auto-generated edit script: updates Leve market-price figures per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 158
$ws.Range("I55").Value = 122.5
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 122.5
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = 91.5
$ws.Range("N55").Value = -728
$ws.Range("H98").Value = 3448.762
$ws.Range("I98").Value = 1789.2941
$ws.Range("J98").Value = 10501.5
$ws.Range("K98").Value = 1789.2941
$ws.Range("L98").Value = 10501.5
$ws.Range("M98").Value = -291.2941000000001
$ws.Range("N98").Value = -13497.5
$ws.Range("H122").Value = 3448.762
$ws.Range("I122").Value = 1789.2941
$ws.Range("J122").Value = 10501.5
$ws.Range("K122").Value = 5367.8823
$ws.Range("L122").Value = 31504.5
$ws.Range("M122").Value = -2917.8823
$ws.Range("N122").Value = -36404.5
$ws.Range("H129").Value = 1013.7647
$ws.Range("J129").Value = 1065.125
$ws.Range("L129").Value = 3195.375
$ws.Range("N129").Value = -13195.375
$ws.Range("H141").Value = 61136.53
$ws.Range("I141").Value = 78725.46000000001
$ws.Range("J141").Value = 3972.5
$ws.Range("K141").Value = 236176.38
$ws.Range("L141").Value = 11917.5
$ws.Range("M141").Value = -230996.38
$ws.Range("N141").Value = -22277.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 32900
$ws.Range("J104").Value = 32900
$ws.Range("L104").Value = 32900
$ws.Range("N104").Value = -39888
$ws.Range("H122").Value = 1792.8148
$ws.Range("I122").Value = 1339.3914
$ws.Range("K122").Value = 4018.1742
$ws.Range("M122").Value = -1568.1742
$ws.Range("H132").Value = 2348.347
$ws.Range("I132").Value = 1825.1892
$ws.Range("K132").Value = 5475.5676
$ws.Range("M132").Value = -2945.5676

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2417.7083
$ws.Range("I134").Value = 1387.1471
$ws.Range("J134").Value = 4920.5
$ws.Range("K134").Value = 4161.4413
$ws.Range("L134").Value = 14761.5
$ws.Range("M134").Value = -1626.4413
$ws.Range("N134").Value = -19831.5
$ws.Range("H140").Value = 54128
$ws.Range("J140").Value = 54128
$ws.Range("L140").Value = 54128
$ws.Range("N140").Value = -64488

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -26232
$ws.Range("H99").Value = 13337177
$ws.Range("I99").Value = 20001864
$ws.Range("K99").Value = 20001864
$ws.Range("M99").Value = -20000366
$ws.Range("H122").Value = 1844.12
$ws.Range("I122").Value = 1254.8889
$ws.Range("J122").Value = 2175.5625
$ws.Range("K122").Value = 3764.6667
$ws.Range("L122").Value = 6526.6875
$ws.Range("M122").Value = -1314.6667
$ws.Range("N122").Value = -11426.6875
$ws.Range("H126").Value = 13337177
$ws.Range("I126").Value = 20001864
$ws.Range("K126").Value = 60005592
$ws.Range("M126").Value = -60003122
$ws.Range("H135").Value = 37384
$ws.Range("J135").Value = 37384
$ws.Range("L135").Value = 37384
$ws.Range("N135").Value = -47524

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7463556
$ws.Range("J131").Value = 906.4603
$ws.Range("L131").Value = 2719.3809
$ws.Range("N131").Value = -12799.3809

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2000
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -2992
$ws.Range("H122").Value = 3116.1
$ws.Range("I122").Value = 1487.6
$ws.Range("K122").Value = 4462.799999999999
$ws.Range("M122").Value = -2012.799999999999
$ws.Range("H126").Value = 2928.99
$ws.Range("I126").Value = 2928.99
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8786.969999999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6316.969999999999
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3132.1614
$ws.Range("I132").Value = 1981.7333
$ws.Range("J132").Value = 4210.6875
$ws.Range("K132").Value = 5945.199900000001
$ws.Range("L132").Value = 12632.0625
$ws.Range("M132").Value = -3415.199900000001
$ws.Range("N132").Value = -17692.0625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4183.0835
$ws.Range("I7").Value = 2240.1
$ws.Range("J7").Value = 5570.9287
$ws.Range("K7").Value = 2240.1
$ws.Range("L7").Value = 5570.9287
$ws.Range("M7").Value = -2128.1
$ws.Range("N7").Value = -5794.9287
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("K30").Value = 1000
$ws.Range("M30").Value = -892
$ws.Range("H32").Value = 9000
$ws.Range("I32").Value = 3000
$ws.Range("K32").Value = 3000
$ws.Range("M32").Value = -2683
$ws.Range("H35").Value = 365.5
$ws.Range("I35").Value = 365.5
$ws.Range("K35").Value = 365.5
$ws.Range("M35").Value = -29.5
$ws.Range("H38").Value = 28000
$ws.Range("J38").Value = 28000
$ws.Range("L38").Value = 28000
$ws.Range("N38").Value = -28820
$ws.Range("H39").Value = 10000
$ws.Range("J39").Value = 10000
$ws.Range("L39").Value = 10000
$ws.Range("N39").Value = -10920
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H51").Value = 42000
$ws.Range("J51").Value = 42000
$ws.Range("L51").Value = 42000
$ws.Range("N51").Value = -42956
$ws.Range("H54").Value = 35083.5
$ws.Range("J54").Value = 35083.5
$ws.Range("L54").Value = 35083.5
$ws.Range("N54").Value = -36371.5
$ws.Range("H56").Value = 12982.6
$ws.Range("I56").Value = 8266.666999999999
$ws.Range("K56").Value = 8266.666999999999
$ws.Range("M56").Value = -7575.666999999999
$ws.Range("H58").Value = 35000
$ws.Range("J58").Value = 35000
$ws.Range("L58").Value = 35000
$ws.Range("N58").Value = -35520
$ws.Range("H122").Value = 4876.5
$ws.Range("I122").Value = 3185.7693
$ws.Range("K122").Value = 9557.3079
$ws.Range("M122").Value = -7107.3079
$ws.Range("H126").Value = 4183.0835
$ws.Range("I126").Value = 2240.1
$ws.Range("J126").Value = 5570.9287
$ws.Range("K126").Value = 6720.299999999999
$ws.Range("L126").Value = 16712.7861
$ws.Range("M126").Value = -4250.299999999999
$ws.Range("N126").Value = -21652.7861
$ws.Range("H132").Value = 3753.426
$ws.Range("I132").Value = 1085.4857
$ws.Range("J132").Value = 8668.053
$ws.Range("K132").Value = 3256.4571
$ws.Range("L132").Value = 26004.159
$ws.Range("M132").Value = -726.4570999999996
$ws.Range("N132").Value = -31064.159

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3554
$ws.Range("I122").Value = 2248.3333
$ws.Range("J122").Value = 5120.8
$ws.Range("K122").Value = 6744.999899999999
$ws.Range("L122").Value = 15362.4
$ws.Range("M122").Value = -4294.999899999999
$ws.Range("N122").Value = -20262.4
$ws.Range("H130").Value = 36131
$ws.Range("J130").Value = 36131
$ws.Range("L130").Value = 36131
$ws.Range("N130").Value = -46171
$ws.Range("H132").Value = 27789536
$ws.Range("I132").Value = 26125
$ws.Range("J132").Value = 41671240
$ws.Range("K132").Value = 78375
$ws.Range("L132").Value = 125013720
$ws.Range("M132").Value = -75845
$ws.Range("N132").Value = -125018780
